$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new column before DT, with a "19-nov" header
# and "-" placeholders for rows 2..25 (a new date column with no data yet) ---
$wsSpot = $wb.Worksheets.Item("Prix Spot")
$wsSpot.Columns("DT").Insert()
$wsSpot.Range("DT1").Value = "19-nov"
for ($r = 2; $r -le 25; $r++) {
    $wsSpot.Cells.Item($r, 124).Value = "-"
}

# --- Sheet "Gaz": append a new row with the next day's price ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A153").Value = "'2025-11-17"
$wsGaz.Range("B153").Value = 30.395

# --- Sheet "CO2": append a new row with the next day's price ---
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A153").Value = "'2025-11-17"
$wsCO2.Range("B153").Value = 79.68000000000001
